$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 48: Intersection of two arrays 2
$ws.Range("A48").Value = 350
$ws.Range("B48").Value = "Intersection of two arrays 2"
$ws.Range("C48").Value = "Dictionary/Sorting/Skipping"

# Add new row 49: Valid Perfect Square
$ws.Range("A49").Value = 367
$ws.Range("B49").Value = "Valid Perfect Square"
$ws.Range("C49").Value = "Binary Search"

# Update the view/selection to match the new state
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D49").Select()
